# [IMP] New data for test environment
#
# Updates a handful of journal codes/names in the "account_journal" sheet
# so they reference the new "external.*" naming scheme instead of the old
# "z0bug.*" one, and nudges a couple of purely cosmetic sheet settings
# (selected cell, column A width) to match the refreshed test workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Journal identifiers / codes -----------------------------------------
# Row 7 : "Operazioni varie" misc journal
$ws.Range("A7").Value = "external.MISC|VARIE"

# Row 8 : "Fatture di vendita" sales journal
$ws.Range("A8").Value = "external.FAT|FATT|INV"
$ws.Range("G8").Value = "FAT"

# Row 9 : "Fatture di acquisto" purchase journal
$ws.Range("A9").Value = "external.ACQ|FATTU|BILL"

# --- Cosmetic sheet view tweaks -------------------------------------------
# Widen column A slightly and move the active selection to A7.
$ws.Columns.Item(1).ColumnWidth = 21.5
$ws.Range("A7").Select() | Out-Null
